$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row heights: rows 2-10 and 13-18 change from 14.4 to 13.8 ---
# (rows 11 and 12 are already 13.8 in the source workbook)
$ws.Range("A2:A10").EntireRow.RowHeight = 13.8
$ws.Range("A13:A18").EntireRow.RowHeight = 13.8

# --- New data: columns V (CO2) and W (O2) get a value of 0 on every
#     data row (2-18); these cells were previously empty/sparse. ---
$ws.Range("V2:W18").Value = 0

# --- View state: scroll so column M is the left-most visible column,
#     and select V2:W18 (the freshly populated block) with V2 active. ---
$ws.Range("V2:W18").Select()
$excel.ActiveWindow.ScrollColumn = 13
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.DisplayGridlines = $true
